$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.499.97"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'3.430.07"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'587.56"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'137.92"
$ws.Range("E6").Value = "  -3.75%  "
$ws.Range("D7").Value = "'3.428.66"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -5.27%  "
$ws.Range("E11").Value = "  -9.25%  "
$ws.Range("E12").Value = "  -7.04%  "
$ws.Range("D13").Value = "'4.010.75"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("E14").Value = "  -9.90%  "
$ws.Range("D15").Value = "'26.35"
$ws.Range("E15").Value = "  -8.40%  "
$ws.Range("D16").Value = "'3.428.25"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "'65.453.66"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "'9.80"
$ws.Range("E19").Value = "  -10.26%  "
$ws.Range("D20").Value = "'5.88"
$ws.Range("E20").Value = "  -5.01%  "
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").Value = "'391.58"
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("D23").Value = "'0.558"
$ws.Range("E23").Value = "  -6.41%  "
$ws.Range("D24").Value = "'73.20"
$ws.Range("E24").Value = "  -5.53%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'3.566.25"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("E27").Value = "  -7.85%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "'7.17"
$ws.Range("E29").Value = "  -7.27%  "
$ws.Range("E30").Value = "  -9.21%  "
$ws.Range("E31").Value = "  -8.53%  "
$ws.Range("D32").Value = "'3.436.22"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'0.145"
$ws.Range("E34").Value = "  -5.76%  "
$ws.Range("D35").Value = "'22.98"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("D36").Value = "'173.00"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").Value = "'6.86"
$ws.Range("E37").Value = "  -8.51%  "
$ws.Range("E38").Value = "  -8.21%  "
$ws.Range("E39").Value = "  -6.69%  "
$ws.Range("E40").Value = "  -8.71%  "
$ws.Range("D41").Value = "'0.0766"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("E42").Value = "  -4.34%  "
$ws.Range("D43").Value = "'43.56"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -12.11%  "
$ws.Range("D46").Value = "'1.62"
$ws.Range("E46").Value = "  -8.81%  "
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D48").Value = "'22.29"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "'6.54"
$ws.Range("E49").Value = "  -7.71%  "
$ws.Range("E50").Value = "  -13.46%  "
$ws.Range("D51").Value = "'2.194.11"
$ws.Range("E51").Value = "  -6.61%  "
